$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.910.84"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").Value = "3.178.27"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.20"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.26"
$ws.Range("E6").Value = "  +5.36%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.178.83"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.01"
$ws.Range("E11").Value = "  +5.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.14"
$ws.Range("E13").Value = "  +5.88%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000252"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "3.713.48"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.121"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.37"
$ws.Range("E17").Value = "  +4.24%  "
$ws.Range("D18").Value = "64.712.38"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "3.189.27"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.76"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.00"
$ws.Range("E21").Value = "  +5.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.765"
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.79"
$ws.Range("E23").Value = "  +4.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.68"
$ws.Range("E24").Value = "  +5.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +10.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.42"
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("E28").Value = "  +8.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.75"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.46"
$ws.Range("E30").Value = "  +6.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.26"
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +8.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.29"
$ws.Range("E34").Value = "  +4.95%  "
$ws.Range("D35").Value = "0.0₃0897"
$ws.Range("E35").Value = "  +5.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.57"
$ws.Range("E36").Value = "  +7.46%  "
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.30"
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "471.83"
$ws.Range("E40").Value = "  +7.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.40"
$ws.Range("E41").Value = "  +6.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.46"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.301"
$ws.Range("E43").Value = "  +8.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0380"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("D45").Value = "2.940.75"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").Value = "  +3.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.69"
$ws.Range("E47").Value = "  +4.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.40"
$ws.Range("E48").Value = "  +4.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.94"
$ws.Range("E49").Value = "  +6.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("E50").Value = "  +7.61%  "
